$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original A:F values for every data row (2-25) before any
# modification, since several rows are swapped/rotated with each other.
$orig = @{}
for ($r = 2; $r -le 25; $r++) {
    $orig[$r] = $ws.Range("A$r`:F$r").Value2
}

# Mapping: new row number -> source row number (values taken from the
# original snapshot above). Rows not listed keep their original values.
$rowMap = @{
    2  = 9
    4  = 14
    5  = 6
    6  = 11
    7  = 12
    8  = 2
    9  = 10
    10 = 4
    11 = 15
    12 = 13
    13 = 7
    14 = 5
    15 = 8
    17 = 20
    20 = 17
    22 = 23
    23 = 22
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $ws.Range("A$destRow`:F$destRow").Value2 = $orig[$srcRow]
}
